# Processed P vs I data: add "alpha" and "pmax" columns (E, F) to each of
# the three data sheets (Rhodymenia, Prionitis, Crypto_control), with a
# per-sheet constant alpha/pmax value repeated down every data row.
# Alpha gets a scientific-notation number format; pmax stays General.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: Rhodymenia ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E1").Value = "alpha"
$ws1.Range("F1").Value = "pmax"
$ws1.Range("E2:E25").Value = 0.000060322529999999999
$ws1.Range("E2:E25").NumberFormat = "0.00E+00"
$ws1.Range("F2:F25").Value = 0.24113119999999999

# ---- Sheet 2: Prionitis ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E1").Value = "alpha"
$ws2.Range("F1").Value = "pmax"
$ws2.Range("E2:E25").Value = 0.00003416555
$ws2.Range("E2:E25").NumberFormat = "0.00E+00"
$ws2.Range("F2:F25").Value = 0.069251270000000004

# ---- Sheet 3: Crypto_control ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("E1").Value = "alpha"
$ws3.Range("F1").Value = "pmax"
$ws3.Range("E2:E25").Value = 0.000056531579999999998
$ws3.Range("E2:E25").NumberFormat = "0.00E+00"
$ws3.Range("F2:F25").Value = 0.21842800000000001

# ---- Restore/update the per-sheet selections (activeCell) ----
$ws1.Activate()
$ws1.Range("H13").Select()

$ws2.Activate()
$ws2.Range("H15").Select()

$ws3.Activate()
$ws3.Range("H9").Select()
